$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.41"
$ws.Range("E2").Value = "'-0.23%"
$ws.Range("G2").Value = "'22"
$ws.Range("D3").Value = "'48.32"
$ws.Range("E3").Value = "'8.30%"
$ws.Range("G3").Value = "'22"
$ws.Range("D4").Value = "'5.244"
$ws.Range("E4").Value = "'1.74%"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.07865"
$ws.Range("E5").Value = "'-2.17%"
$ws.Range("G5").Value = "'22"
$ws.Range("D6").Value = "'4.576"
$ws.Range("E6").Value = "'1.29%"
$ws.Range("G6").Value = "'22"
$ws.Range("D7").Value = "'1.326"
$ws.Range("E7").Value = "'22.34%"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'1.601"
$ws.Range("E8").Value = "'-3.34%"
$ws.Range("G8").Value = "'22"
$ws.Range("D9").Value = "'0.1240"
$ws.Range("E9").Value = "'-4.82%"
$ws.Range("G9").Value = "'22"
$ws.Range("D10").Value = "'0.1944"
$ws.Range("E10").Value = "'1.37%"
$ws.Range("G10").Value = "'22"
$ws.Range("D11").Value = "'0.09391"
$ws.Range("E11").Value = "'0.09%"
$ws.Range("G11").Value = "'22"
$ws.Range("D12").Value = "'0.04538"
$ws.Range("E12").Value = "'7.59%"
$ws.Range("G12").Value = "'22"
$ws.Range("D13").Value = "'0.1046"
$ws.Range("E13").Value = "'0.65%"
$ws.Range("G13").Value = "'22"
$ws.Range("D14").Value = "'0.001307"
$ws.Range("E14").Value = "'-0.49%"
$ws.Range("G14").Value = "'22"
$ws.Range("D15").Value = "'0.04199"
$ws.Range("E15").Value = "'-0.24%"
$ws.Range("G15").Value = "'22"
$ws.Range("D16").Value = "'0.005825"
$ws.Range("E16").Value = "'-1.18%"
$ws.Range("G16").Value = "'22"
$ws.Range("E17").Value = "'-1.32%"
$ws.Range("G17").Value = "'22"
$ws.Range("D18").Value = "'2.467"
$ws.Range("E18").Value = "'2.72%"
$ws.Range("G18").Value = "'22"
$ws.Range("E19").Value = "'2.01%"
$ws.Range("G19").Value = "'22"
$ws.Range("D20").Value = "'8.082"
$ws.Range("E20").Value = "'1.12%"
$ws.Range("G20").Value = "'22"
$ws.Range("E21").Value = "'-0.44%"
$ws.Range("G21").Value = "'22"
$ws.Range("D22").Value = "'0.3091"
$ws.Range("E22").Value = "'-1.48%"
$ws.Range("G22").Value = "'22"
$ws.Range("E23").Value = "'1.55%"
$ws.Range("G23").Value = "'22"
$ws.Range("D24").Value = "'0.004224"
$ws.Range("G24").Value = "'22"
$ws.Range("D25").Value = "'0.0001356"
$ws.Range("E25").Value = "'1.30%"
$ws.Range("G25").Value = "'22"
$ws.Range("D26").Value = "'0.0003548"
$ws.Range("G26").Value = "'22"
$ws.Range("G27").Value = "'22"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("D38").Value = "'0.02624"
$ws.Range("E38").Value = "'-1.83%"
$ws.Range("G38").Value = "'22"
$ws.Range("D39").Value = "'0.05853"
$ws.Range("E39").Value = "'8.17%"
$ws.Range("G39").Value = "'22"
$ws.Range("D40").Value = "'0.01081"
$ws.Range("E40").Value = "'92.10%"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.007997"
$ws.Range("E41").Value = "'3.13%"
$ws.Range("G41").Value = "'22"
$ws.Range("D42").Value = "'0.1440"
$ws.Range("E42").Value = "'1.64%"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.008325"
$ws.Range("E43").Value = "'13.41%"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.008562"
$ws.Range("E44").Value = "'8.05%"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.3134"
$ws.Range("E45").Value = "'0.28%"
$ws.Range("G45").Value = "'22"
$ws.Range("D46").Value = "'0.00006935"
$ws.Range("E46").Value = "'2.03%"
$ws.Range("G46").Value = "'22"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.90%"
$ws.Range("G47").Value = "'22"
$ws.Range("E48").Value = "'-6.82%"
$ws.Range("G48").Value = "'22"
$ws.Range("D49").Value = "'0.004010"
$ws.Range("E49").Value = "'0.93%"
$ws.Range("G49").Value = "'22"
$ws.Range("D50").Value = "'0.00002105"
$ws.Range("E50").Value = "'0.90%"
$ws.Range("G50").Value = "'22"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("E51").Value = "'0.90%"
$ws.Range("G51").Value = "'22"
